$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'262.34"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'22.95"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'6.184"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.06099"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'6.736"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'3.458"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").Value = "'0.7965"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'0.1588"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.08082"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'0.03527"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'0.03081"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'0.09316"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'3.849"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'0.001707"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'0.04785"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").Value = "'0.006194"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'0.001092"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'0.003680"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Value = "'3.705"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "'2.223"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'0.3364"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Value = "'0.1253"
$ws.Range("D26").Style = "Normal"
$ws.Range("D40").Value = "'0.04603"
$ws.Range("D40").Style = "Normal"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.007127"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("D42").Value = "'0.003899"
$ws.Range("D42").Style = "Normal"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").Value = "'0.1118"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "42BKEXTokenBKK"
$ws.Range("D44").Value = "'0.01065"
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").Value = "'0.00005916"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").Value = "'0.6997"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Value = "'0.09178"
$ws.Range("D49").Style = "Normal"
